$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) "Data" sheet: refresh the whole Fecha/Valor series.
#    The series now starts two years later (2023, 2022 added on
#    top) and extends two years further back (1991, 1990 appended
#    at the bottom), so the sheet grows from 32 to 34 data rows.
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Data")

# Make room for the two newest years by inserting two fresh rows
# right below the header (row 1); everything else shifts down.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

$years = @("2023","2022","2021","2020","2019","2018","2017","2016","2015","2014","2013","2012","2011","2010","2009","2008","2007","2006","2005","2004","2003","2002","2001","2000","1999","1998","1997","1996","1995","1994","1993","1992","1991","1990")
$values = @(905.8,871.3,821.1,831.5,874.5,846.5,826.1,780.9,738.3,747.4,716.6,678.9,638.3,524,510.9,447.5,380.8,341.1,317,295.3,299.5,292.2,318.3,302.2,313.4,188.7,183.8,173.2,171.6,186.1,167.8,144.2,131,133.2)

# Keep the "Fecha" column as text (years would otherwise be
# auto-coerced to numbers) before writing the values.
$ws.Range("A2:A35").NumberFormat = "@"

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $years[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# ---------------------------------------------------------------
# 2) "Metadata" sheet: add an "actualizacion" / "Julio 2025" row
#    right before the "cita" row. The blank A1 header cell is also
#    normalised to match B1 (a single space) in the refreshed file.
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Metadata")
$ws2.Range("A1").Value = " "
$ws2.Rows.Item(9).Insert()
$ws2.Range("A9").Value = "actualizacion"
$ws2.Range("B9").Value = "Julio 2025"
